$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-10 (Q0..Q8) with new values for columns B..G
$ws.Range("B2").Value = 0.1340334775780207
$ws.Range("C2").Value = 0.51771755592959
$ws.Range("D2").Value = 0.6327526057370099
$ws.Range("E2").Value = 0.7954574820422584
$ws.Range("F2").Value = 0.7918859673449575
$ws.Range("G2").Value = 51

$ws.Range("B3").Value = -0.04416780937479916
$ws.Range("C3").Value = 0.5817619418406447
$ws.Range("D3").Value = 0.71848441568627
$ws.Range("E3").Value = 0.8476346003357048
$ws.Range("F3").Value = 0.8550770464756325
$ws.Range("G3").Value = 50

$ws.Range("B4").Value = 0.1307748531607076
$ws.Range("C4").Value = 0.567321327886076
$ws.Range("D4").Value = 0.6923438032425172
$ws.Range("E4").Value = 0.8320719940260681
$ws.Range("F4").Value = 0.8302465159786168
$ws.Range("G4").Value = 49

$ws.Range("B5").Value = -0.05427852441735345
$ws.Range("C5").Value = 0.5770490505655083
$ws.Range("D5").Value = 0.6815373014256153
$ws.Range("E5").Value = 0.8255527248005516
$ws.Range("F5").Value = 0.8324837852131098
$ws.Range("G5").Value = 48

$ws.Range("B6").Value = 0.09860079171026061
$ws.Range("C6").Value = 0.5763606548324329
$ws.Range("D6").Value = 0.6635861464005488
$ws.Range("E6").Value = 0.8146079709900639
$ws.Range("F6").Value = 0.8173606705215342
$ws.Range("G6").Value = 47

$ws.Range("B7").Value = -0.09963899298847101
$ws.Range("C7").Value = 0.5177403320122128
$ws.Range("D7").Value = 0.5601913235939101
$ws.Range("E7").Value = 0.7484592998913903
$ws.Range("F7").Value = 0.7499943133832871
$ws.Range("G7").Value = 46

$ws.Range("B8").Value = 0.04033191057269756
$ws.Range("C8").Value = 0.5192072963624517
$ws.Range("D8").Value = 0.5728721284028437
$ws.Range("E8").Value = 0.7568831669437786
$ws.Range("F8").Value = 0.764348295529329
$ws.Range("G8").Value = 45

$ws.Range("B9").Value = -0.05085408857912471
$ws.Range("C9").Value = 0.4865973670774532
$ws.Range("D9").Value = 0.5115951227155482
$ws.Range("E9").Value = 0.71525878024359
$ws.Range("F9").Value = 0.721696891105934
$ws.Range("G9").Value = 44

$ws.Range("B10").Value = 0.0175971372297503
$ws.Range("C10").Value = 0.5253937883667443
$ws.Range("D10").Value = 0.5760576781307024
$ws.Range("E10").Value = 0.7589846362942417
$ws.Range("F10").Value = 0.7677605779513007
$ws.Range("G10").Value = 43

# Add new row 11 (Q9) with the same style as the other label cells
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Q9"

$ws.Range("B11").Value = -0.03005877115108942
$ws.Range("C11").Value = 0.4860398704780853
$ws.Range("D11").Value = 0.4986732487203297
$ws.Range("E11").Value = 0.7061680031836119
$ws.Range("F11").Value = 0.7140801382553689
$ws.Range("G11").Value = 42
